$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate rows 2-8 into a single Python-tuple-like string in A2
$ws.Range("A2").Value = '(''Garruk Wildspeaker'', [''{2}{G}{G}'', ''Legendary Planeswalker ' + [char]0x2014 + ' Garruk'', ''+1: Untap two target lands.'', ''' + [char]0x2212 + '1: Create a 3/3 green Beast creature token.'', ''' + [char]0x2212 + '4: Creatures you control get +3/+3 and gain trample until end of turn.'', ''Loyalty: 3''])'

# Remove the now-unused rows 3-8 that used to hold the individual fields
$ws.Rows("3:8").Delete()
